# Update final evaluation results across the three result sheets.
$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B2").Value = 0.5374531835205992
$wsSummary.Range("C2").Value = 0.519455252918288
$wsSummary.Range("D2").Value = 1
$wsSummary.Range("E2").Value = 0.6837387964148528
$wsSummary.Range("F2").Value = 0.8438685208596713
$wsSummary.Range("G2").Value = 0.9656419529837251
$wsSummary.Range("H2").Value = 0.7046108095218055
$wsSummary.Range("I2").Value = 534
$wsSummary.Range("J2").Value = 494
$wsSummary.Range("K2").Value = 40
$wsSummary.Range("L2").Value = 0

# --- Classification Report sheet ---
$wsClass = $wb.Worksheets.Item("Classification Report")
$wsClass.Range("B2").Value = 1
$wsClass.Range("C2").Value = 0.0749063670411985
$wsClass.Range("D2").Value = 0.1393728222996516

$wsClass.Range("B3").Value = 0.519455252918288
$wsClass.Range("C3").Value = 1
$wsClass.Range("D3").Value = 0.6837387964148528

$wsClass.Range("B4").Value = 0.5374531835205992
$wsClass.Range("C4").Value = 0.5374531835205992
$wsClass.Range("D4").Value = 0.5374531835205992
$wsClass.Range("E4").Value = 0.5374531835205992

$wsClass.Range("B5").Value = 0.7597276264591439
$wsClass.Range("C5").Value = 0.5374531835205992
$wsClass.Range("D5").Value = 0.4115558093572522

$wsClass.Range("B6").Value = 0.7597276264591439
$wsClass.Range("C6").Value = 0.5374531835205992
$wsClass.Range("D6").Value = 0.4115558093572522

# --- Confusion Matrix sheet ---
$wsConf = $wb.Worksheets.Item("Confusion Matrix")
$wsConf.Range("B2").Value = 40
$wsConf.Range("C2").Value = 494

$wsConf.Range("B3").Value = 0
$wsConf.Range("C3").Value = 534
